$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 271, shifting rows 271:295 down to 272:296
$ws.Rows.Item(271).Insert()

# Populate the newly inserted row 271 with its new data
$ws.Cells.Item(271, 1).Value = 10
$ws.Cells.Item(271, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(271, 3).Value = "La Araucanía"
$ws.Cells.Item(271, 4).Value = 44461
$ws.Cells.Item(271, 5).Value = 9
$ws.Cells.Item(271, 6).Value = 100112043
$ws.Cells.Item(271, 7).Value = "Pepino ensalada"
$ws.Cells.Item(271, 8).Value = "Sin especificar"
$ws.Cells.Item(271, 9).Value = "Primera"
$ws.Cells.Item(271, 10).Value = 140
$ws.Cells.Item(271, 11).Value = 17000
$ws.Cells.Item(271, 12).Value = 18000
$ws.Cells.Item(271, 13).Value = 17714
$ws.Cells.Item(271, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(271, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(271, 16).Value = 295
$ws.Cells.Item(271, 17).Value = 60
$ws.Cells.Item(271, 18).Value = "Hortaliza"
